$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 with same style as the other header cells (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Add the time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 10:50:43.895758"
$ws.Range("F3").Value = "2021-10-05 10:50:43.895768"
$ws.Range("F4").Value = "2021-10-05 10:50:43.895772"
